$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the very top and fill it with the
# "Chức năng" / "Mô tả" column headers.
$ws.Rows.Item(1).Insert()
$ws.Cells.Item(1, 1).Value = "Chức năng"
$ws.Cells.Item(1, 2).Value = "Mô tả"

# Fix the typo in the "2.1, Thêm điện thoại" item (comma -> period); after
# the header-row insert above it now lives on row 6.
$ws.Cells.Item(6, 1).Value = "2.1. Thêm điện thoại"

# Give the new description column (B) a sensible width, matching column A's
# styling convention (~50.71 characters wide).
$ws.Columns.Item(2).ColumnWidth = 49.8

# Restore a sane view: cursor on B8.
$ws.Range("B8").Select()
